$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update swapped / rotated match rows (columns F:V) ----

# Row 13
$ws.Range("F13").Value = 'Latina'
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 'Potenza'
$ws.Range("I13").Value = 1
$ws.Range("J13").Value = 2.27
$ws.Range("K13").Value = '07/09/2023 15:12'
$ws.Range("L13").Value = 2.74
$ws.Range("M13").Value = '10/09/2023 20:44'
$ws.Range("N13").Value = 2.91
$ws.Range("O13").Value = '07/09/2023 15:12'
$ws.Range("P13").Value = 2.87
$ws.Range("Q13").Value = '10/09/2023 20:30'
$ws.Range("R13").Value = 3.24
$ws.Range("S13").Value = '07/09/2023 15:12'
$ws.Range("T13").Value = 2.91
$ws.Range("U13").Value = '10/09/2023 20:44'
$ws.Range("V13").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/latina-potenza/8IavFPbJ/'

# Row 14
$ws.Range("F14").Value = 'Monopoli'
$ws.Range("G14").Value = 2
$ws.Range("H14").Value = 'Monterosi'
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 1.86
$ws.Range("K14").Value = '07/09/2023 15:12'
$ws.Range("L14").Value = 1.99
$ws.Range("M14").Value = '10/09/2023 20:37'
$ws.Range("N14").Value = 3.09
$ws.Range("O14").Value = '07/09/2023 15:12'
$ws.Range("P14").Value = 3.05
$ws.Range("Q14").Value = '10/09/2023 20:36'
$ws.Range("R14").Value = 4.14
$ws.Range("S14").Value = '07/09/2023 15:12'
$ws.Range("T14").Value = 4.41
$ws.Range("U14").Value = '10/09/2023 20:37'
$ws.Range("V14").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-monterosi/CjQKNH2J/'

# Row 23
$ws.Range("F23").Value = 'Avellino'
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 'Foggia'
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 2.21
$ws.Range("K23").Value = '15/09/2023 04:42'
$ws.Range("L23").Value = 2.18
$ws.Range("M23").Value = '17/09/2023 20:38'
$ws.Range("N23").Value = 2.91
$ws.Range("O23").Value = '15/09/2023 04:42'
$ws.Range("P23").Value = 3
$ws.Range("Q23").Value = '17/09/2023 20:38'
$ws.Range("R23").Value = 3.38
$ws.Range("S23").Value = '15/09/2023 04:42'
$ws.Range("T23").Value = 3.76
$ws.Range("U23").Value = '17/09/2023 20:38'
$ws.Range("V23").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/avellino-foggia/OrPbGZfg/'

# Row 24
$ws.Range("F24").Value = 'Casertana'
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 'Benevento'
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 3.61
$ws.Range("K24").Value = '17/09/2023 08:52'
$ws.Range("L24").Value = 3.72
$ws.Range("M24").Value = '17/09/2023 20:34'
$ws.Range("N24").Value = 3.24
$ws.Range("O24").Value = '17/09/2023 08:52'
$ws.Range("P24").Value = 3.33
$ws.Range("Q24").Value = '17/09/2023 20:34'
$ws.Range("R24").Value = 2.09
$ws.Range("S24").Value = '17/09/2023 08:52'
$ws.Range("T24").Value = 2.05
$ws.Range("U24").Value = '17/09/2023 20:34'
$ws.Range("V24").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/casertana-benevento/0ATjgTNq/'

# Row 72
$ws.Range("F72").Value = 'Juve Stabia'
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 'Catania'
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2.66
$ws.Range("K72").Value = '12/10/2023 08:13'
$ws.Range("L72").Value = 2.82
$ws.Range("M72").Value = '15/10/2023 13:51'
$ws.Range("N72").Value = 2.83
$ws.Range("O72").Value = '12/10/2023 08:13'
$ws.Range("P72").Value = 2.99
$ws.Range("Q72").Value = '15/10/2023 13:51'
$ws.Range("R72").Value = 2.69
$ws.Range("S72").Value = '12/10/2023 08:13'
$ws.Range("T72").Value = 2.72
$ws.Range("U72").Value = '15/10/2023 13:51'
$ws.Range("V72").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/juve-stabia-catania/pvJB7aDE/'

# Row 73
$ws.Range("F73").Value = 'ACR Messina'
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 'Giugliano'
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 2.14
$ws.Range("K73").Value = '12/10/2023 08:13'
$ws.Range("L73").Value = 2.12
$ws.Range("M73").Value = '15/10/2023 13:51'
$ws.Range("N73").Value = 3
$ws.Range("O73").Value = '12/10/2023 08:13'
$ws.Range("P73").Value = 3.11
$ws.Range("Q73").Value = '15/10/2023 13:51'
$ws.Range("R73").Value = 3.32
$ws.Range("S73").Value = '12/10/2023 08:13'
$ws.Range("T73").Value = 3.79
$ws.Range("U73").Value = '15/10/2023 13:51'
$ws.Range("V73").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/acr-messina-giugliano/vsSLpx61/'

# Row 74
$ws.Range("F74").Value = 'Benevento'
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 'Picerno'
$ws.Range("I74").Value = 2
$ws.Range("J74").Value = 1.71
$ws.Range("K74").Value = '12/10/2023 08:13'
$ws.Range("L74").Value = 1.68
$ws.Range("M74").Value = '15/10/2023 13:55'
$ws.Range("N74").Value = 3.33
$ws.Range("O74").Value = '12/10/2023 08:13'
$ws.Range("P74").Value = 3.55
$ws.Range("Q74").Value = '15/10/2023 13:55'
$ws.Range("R74").Value = 4.59
$ws.Range("S74").Value = '12/10/2023 08:13'
$ws.Range("T74").Value = 5.52
$ws.Range("U74").Value = '15/10/2023 13:55'
$ws.Range("V74").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/benevento-picerno/hpWPqdL7/'

# Row 94
$ws.Range("F94").Value = 'Casertana'
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 'Juve Stabia'
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 2.56
$ws.Range("K94").Value = '24/10/2023 12:42'
$ws.Range("L94").Value = 2.69
$ws.Range("M94").Value = '25/10/2023 20:37'
$ws.Range("N94").Value = 2.83
$ws.Range("O94").Value = '24/10/2023 12:42'
$ws.Range("P94").Value = 3.01
$ws.Range("Q94").Value = '25/10/2023 20:37'
$ws.Range("R94").Value = 2.8
$ws.Range("S94").Value = '24/10/2023 12:42'
$ws.Range("T94").Value = 2.83
$ws.Range("U94").Value = '25/10/2023 20:37'
$ws.Range("V94").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/casertana-juve-stabia/vL2bC8UL/'

# Row 95
$ws.Range("F95").Value = 'Monopoli'
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = 'Picerno'
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = 2.59
$ws.Range("K95").Value = '24/10/2023 12:42'
$ws.Range("L95").Value = 2.67
$ws.Range("M95").Value = '25/10/2023 20:41'
$ws.Range("N95").Value = 2.92
$ws.Range("O95").Value = '24/10/2023 12:42'
$ws.Range("P95").Value = 2.97
$ws.Range("Q95").Value = '25/10/2023 20:38'
$ws.Range("R95").Value = 2.69
$ws.Range("S95").Value = '24/10/2023 12:42'
$ws.Range("T95").Value = 2.9
$ws.Range("U95").Value = '25/10/2023 20:41'
$ws.Range("V95").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-picerno/jD1hGTen/'

# Row 96
$ws.Range("F96").Value = 'Taranto'
$ws.Range("G96").Value = 3
$ws.Range("H96").Value = 'Turris'
$ws.Range("I96").Value = 1
$ws.Range("J96").Value = 2.2
$ws.Range("K96").Value = '24/10/2023 12:42'
$ws.Range("L96").Value = 1.79
$ws.Range("M96").Value = '25/10/2023 20:44'
$ws.Range("N96").Value = 3.04
$ws.Range("O96").Value = '24/10/2023 12:42'
$ws.Range("P96").Value = 3.16
$ws.Range("Q96").Value = '25/10/2023 20:41'
$ws.Range("R96").Value = 3.24
$ws.Range("S96").Value = '24/10/2023 12:42'
$ws.Range("T96").Value = 5.42
$ws.Range("U96").Value = '25/10/2023 20:41'
$ws.Range("V96").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/taranto-turris/vZg5DRu5/'

# Row 107
$ws.Range("F107").Value = 'Virtus Francavilla'
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 'Taranto'
$ws.Range("I107").Value = 2
$ws.Range("J107").Value = 2.15
$ws.Range("K107").Value = '27/10/2023 02:42'
$ws.Range("L107").Value = 2.78
$ws.Range("M107").Value = '29/10/2023 20:41'
$ws.Range("N107").Value = 2.95
$ws.Range("O107").Value = '27/10/2023 02:42'
$ws.Range("P107").Value = 2.82
$ws.Range("Q107").Value = '29/10/2023 20:41'
$ws.Range("R107").Value = 3.36
$ws.Range("S107").Value = '27/10/2023 02:42'
$ws.Range("T107").Value = 2.91
$ws.Range("U107").Value = '29/10/2023 20:41'
$ws.Range("V107").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/virtus-francavilla-taranto/bcTsRpA4/'

# Row 108
$ws.Range("F108").Value = 'Brindisi'
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 'Monopoli'
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 2.15
$ws.Range("K108").Value = '27/10/2023 02:42'
$ws.Range("L108").Value = 2.32
$ws.Range("M108").Value = '29/10/2023 20:41'
$ws.Range("N108").Value = 2.95
$ws.Range("O108").Value = '27/10/2023 02:42'
$ws.Range("P108").Value = 3.04
$ws.Range("Q108").Value = '29/10/2023 20:41'
$ws.Range("R108").Value = 3.36
$ws.Range("S108").Value = '27/10/2023 02:42'
$ws.Range("T108").Value = 3.36
$ws.Range("U108").Value = '29/10/2023 20:41'
$ws.Range("V108").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/brindisi-monopoli/z1mEB5AH/'

# Row 113
$ws.Range("F113").Value = 'Avellino'
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 'Virtus Francavilla'
$ws.Range("I113").Value = 1
$ws.Range("J113").Value = 1.39
$ws.Range("K113").Value = '02/11/2023 08:12'
$ws.Range("L113").Value = 1.41
$ws.Range("M113").Value = '04/11/2023 18:29'
$ws.Range("N113").Value = 3.87
$ws.Range("O113").Value = '02/11/2023 08:12'
$ws.Range("P113").Value = 4.06
$ws.Range("Q113").Value = '04/11/2023 18:29'
$ws.Range("R113").Value = 8.31
$ws.Range("S113").Value = '02/11/2023 08:12'
$ws.Range("T113").Value = 10.1
$ws.Range("U113").Value = '04/11/2023 18:29'
$ws.Range("V113").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/avellino-virtus-francavilla/fPMjPOvH/'

# Row 114
$ws.Range("F114").Value = 'Latina'
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 'Crotone'
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 2.79
$ws.Range("K114").Value = '03/11/2023 13:42'
$ws.Range("L114").Value = 3.12
$ws.Range("M114").Value = '04/11/2023 18:26'
$ws.Range("N114").Value = 2.79
$ws.Range("O114").Value = '03/11/2023 13:42'
$ws.Range("P114").Value = 3.24
$ws.Range("Q114").Value = '04/11/2023 18:26'
$ws.Range("R114").Value = 2.61
$ws.Range("S114").Value = '03/11/2023 13:42'
$ws.Range("T114").Value = 2.34
$ws.Range("U114").Value = '04/11/2023 18:26'
$ws.Range("V114").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/latina-crotone/8W2RIMHp/'

# Row 123
$ws.Range("F123").Value = 'Crotone'
$ws.Range("G123").Value = 2
$ws.Range("H123").Value = 'Monterosi'
$ws.Range("I123").Value = 1
$ws.Range("J123").Value = 1.33
$ws.Range("K123").Value = '09/11/2023 09:13'
$ws.Range("L123").Value = 1.34
$ws.Range("M123").Value = '09/11/2023 14:58'
$ws.Range("N123").Value = 4.71
$ws.Range("O123").Value = '09/11/2023 09:13'
$ws.Range("P123").Value = 5.03
$ws.Range("Q123").Value = '12/11/2023 12:59'
$ws.Range("R123").Value = 8.039999999999999
$ws.Range("S123").Value = '09/11/2023 09:13'
$ws.Range("T123").Value = 9.02
$ws.Range("U123").Value = '12/11/2023 12:15'
$ws.Range("V123").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/crotone-monterosi/zs36Jg8e/'

# Row 124
$ws.Range("F124").Value = 'Turris'
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 'Monopoli'
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = 2.16
$ws.Range("K124").Value = '09/11/2023 09:13'
$ws.Range("L124").Value = 2.33
$ws.Range("M124").Value = '12/11/2023 13:46'
$ws.Range("N124").Value = 2.9
$ws.Range("O124").Value = '09/11/2023 09:13'
$ws.Range("P124").Value = 3.35
$ws.Range("Q124").Value = '12/11/2023 13:46'
$ws.Range("R124").Value = 3.41
$ws.Range("S124").Value = '09/11/2023 09:13'
$ws.Range("T124").Value = 3.03
$ws.Range("U124").Value = '12/11/2023 13:42'
$ws.Range("V124").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/turris-monopoli/l6hiPqdE/'

# Row 138
$ws.Range("F138").Value = 'Avellino'
$ws.Range("G138").Value = 1
$ws.Range("H138").Value = 'Giugliano'
$ws.Range("I138").Value = 3
$ws.Range("J138").Value = 1.4
$ws.Range("K138").Value = '16/11/2023 09:12'
$ws.Range("L138").Value = 1.3
$ws.Range("M138").Value = '19/11/2023 16:11'
$ws.Range("N138").Value = 4.02
$ws.Range("O138").Value = '16/11/2023 09:12'
$ws.Range("P138").Value = 5.21
$ws.Range("Q138").Value = '19/11/2023 16:11'
$ws.Range("R138").Value = 7.15
$ws.Range("S138").Value = '16/11/2023 09:12'
$ws.Range("T138").Value = 10.42
$ws.Range("U138").Value = '19/11/2023 16:11'
$ws.Range("V138").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/avellino-giugliano/rTieO3BK/'

# Row 139
$ws.Range("F139").Value = 'Monopoli'
$ws.Range("G139").Value = 3
$ws.Range("H139").Value = 'Benevento'
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 3.58
$ws.Range("K139").Value = '16/11/2023 09:12'
$ws.Range("L139").Value = 3.22
$ws.Range("M139").Value = '19/11/2023 16:10'
$ws.Range("N139").Value = 3.09
$ws.Range("O139").Value = '16/11/2023 09:12'
$ws.Range("P139").Value = 3.14
$ws.Range("Q139").Value = '19/11/2023 16:14'
$ws.Range("R139").Value = 2.01
$ws.Range("S139").Value = '16/11/2023 09:12'
$ws.Range("T139").Value = 2.33
$ws.Range("U139").Value = '19/11/2023 16:11'
$ws.Range("V139").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-benevento/lUzvGuJ0/'

# Row 140
$ws.Range("F140").Value = 'Potenza'
$ws.Range("G140").Value = 2
$ws.Range("H140").Value = 'Audace Cerignola'
$ws.Range("I140").Value = 2
$ws.Range("J140").Value = 2.29
$ws.Range("K140").Value = '16/11/2023 09:12'
$ws.Range("L140").Value = 2.72
$ws.Range("M140").Value = '19/11/2023 18:26'
$ws.Range("N140").Value = 2.94
$ws.Range("O140").Value = '16/11/2023 09:12'
$ws.Range("P140").Value = 3.11
$ws.Range("Q140").Value = '19/11/2023 18:26'
$ws.Range("R140").Value = 3.07
$ws.Range("S140").Value = '16/11/2023 09:12'
$ws.Range("T140").Value = 2.71
$ws.Range("U140").Value = '19/11/2023 18:26'
$ws.Range("V140").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/potenza-audace-cerignola/UwqnEJmD/'

# Row 141
$ws.Range("F141").Value = 'Juve Stabia'
$ws.Range("G141").Value = 0
$ws.Range("H141").Value = 'Sorrento'
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 1.61
$ws.Range("K141").Value = '16/11/2023 09:12'
$ws.Range("L141").Value = 1.5
$ws.Range("M141").Value = '19/11/2023 18:28'
$ws.Range("N141").Value = 3.49
$ws.Range("O141").Value = '16/11/2023 09:12'
$ws.Range("P141").Value = 3.78
$ws.Range("Q141").Value = '19/11/2023 18:28'
$ws.Range("R141").Value = 5.18
$ws.Range("S141").Value = '16/11/2023 09:12'
$ws.Range("T141").Value = 8.039999999999999
$ws.Range("U141").Value = '19/11/2023 18:28'
$ws.Range("V141").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/juve-stabia-sorrento/t4XVH1ll/'

# ---- Append new match row 142 ----

$ws.Range("A141").Copy($ws.Range("A142"))
$ws.Range("A142").Value = 141
$ws.Range("B142").Value = 'italy'
$ws.Range("C142").Value = 'serie-c-group-c'
$ws.Range("D142").Value = '2023-2024'
$ws.Range("E142").NumberFormat = $ws.Range("E141").NumberFormat
$ws.Range("E142").Value = 45254.86458333334
$ws.Range("F142").Value = 'Foggia'
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 'Latina'
$ws.Range("I142").Value = 2
$ws.Range("J142").Value = 2.03
$ws.Range("K142").Value = '22/11/2023 21:12'
$ws.Range("L142").Value = 2.05
$ws.Range("M142").Value = '24/11/2023 20:41'
$ws.Range("N142").Value = 3
$ws.Range("O142").Value = '22/11/2023 21:12'
$ws.Range("P142").Value = 3.01
$ws.Range("Q142").Value = '24/11/2023 20:41'
$ws.Range("R142").Value = 3.62
$ws.Range("S142").Value = '22/11/2023 21:12'
$ws.Range("T142").Value = 4.23
$ws.Range("U142").Value = '24/11/2023 20:41'
$ws.Range("V142").Value = 'https://www.betexplorer.com/football/italy/serie-c-group-c/foggia-latina/zuC8VKIm/'
